# Update the coffee riskDriverType PnL factor calibrations with the
# re-estimated values (std moving average based re-calibration).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.0006435182975611956
$ws.Range("B3").Value = 0.7821501015054658
$ws.Range("B4").Value = 0.1033347820053903
$ws.Range("B5").Value = "[1.0, 0.022207096170267747, -0.03823237276396907, -0.09617065369444644, -0.10263211315825739, 0.0758623494219552, 0.019710505444981916, 0.04453455219793419, 0.035113853860857094, -0.014984468043391866, -0.025557892622576375, -0.026438729512850374, -0.004154538204781508, 0.004932830434421905, 0.055254612870286414, 0.019036929854697612, 0.01126241651045639, -0.03866852993032194, -0.049356888036541535, -0.0010941808226047954]"

$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1866208303222265
$ws.Range("B5").Value = 0.5700724211969826
$ws.Range("B6").Value = 0.06631996706113626
$ws.Range("B7").Value = 0.1680596614087125
$ws.Range("B8").Value = 0.5811439645546888
$ws.Range("B9").Value = 0.05718170724591762
$ws.Range("B10").Value = "[1.0, 0.09942809580729636, -0.01581844129143008, -0.016360545666903054, 0.016260304349724954, 0.03617349924709207, -0.0016137658443848362, 0.011835130770821649, 0.003275445819416438, 0.008316791367869447, -0.00714019547087173, 0.006778493546845363, -0.016662795757746908, -0.002734101564051823, 0.01009811336756387, -0.0021922430889197322, -0.00834100108288927, 0.005493944931717108, -0.010409530770728885, 0.008474816117425002]"

$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.0007205375476518765
$ws.Range("B3").Value = 0.1113320037510603
$ws.Range("B4").Value = 0.03683693758144928
$ws.Range("B5").Value = [double]"2.806299768179604E-16"
$ws.Range("B6").Value = "[1.0, -0.007700235263675449, -0.04166627016872811, -0.09392418823806951, -0.11406974978898479, 0.10703502772690926, 0.004600835732978159, 0.0555857369761117, 0.027245472303308156, -0.012835634061624553, -0.02741135623324103, -0.03456039081038182, -0.016613775634588435, 0.01971567750553987, 0.047940690903815655, 0.011754148990867191, 0.0010339794633052084, -0.03692751952828282, -0.05183900332966369, -0.001382757226537816]"

$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.001747518788148759
$ws.Range("B3").Value = 0.1120460552106263
$ws.Range("B4").Value = [double]"8.98463376722993E-12"
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[1.0, -0.0037175150995300086, -0.03975663958297556, -0.09489494262139268, -0.11354582316914354, 0.10644449999531425, 0.006231207634274265, 0.05543860190984388, 0.028091415065847537, -0.010960779771878329, -0.027133972791784082, -0.03415335634796045, -0.017068038154701388, 0.01890522966725043, 0.04717232090924073, 0.01269366310501702, 0.001073059811089547, -0.03708449431029641, -0.05219276583871505, -0.00041264933388403903]"
$ws.Range("B7").Value = 0.05818298852740635

$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.0006595151091066017
$ws.Range("B3").Value = 0.1002023654162825
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = [double]"6.594045302147165E-09"
$ws.Range("B6").Value = "[1.0, -0.010922272246122081, -0.038043610632496366, -0.09465962922924305, -0.10367939252057189, 0.07765434801432926, 0.019060563107717327, 0.04410206164273828, 0.03677540085070315, -0.012005608770991008, -0.024799541670654835, -0.026836059697918066, -0.002826758053029191, 0.0038152630475419896, 0.05651319682354552, 0.018403705599285697, 0.012535391725134182, -0.038380745391197284, -0.05021377939301723, -0.003083215758436748]"
$ws.Range("B7").Value = 0.05469225589322876
$ws.Range("B9").Value = 0.7830964463375857

